# Weekly update: insert 3 new rows of data (most recent week) above the
# existing row 367, pushing the old rows 367-379 down to 370-382.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 367 (shifts existing 367.. down by 3)
$ws.Rows.Item(367).Insert()
$ws.Rows.Item(367).Insert()
$ws.Rows.Item(367).Insert()

# Common values shared across these rows
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$catId     = 100112023
$categoria = "Brócoli"
$variedad  = "Sin especificar"
$unidadCom = "`$/unidad"
$origen    = "Región de Arica y Parinacota"
$kgUnid    = 1
$tipo      = "Hortaliza"

function Set-Row($r, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $catId
    $ws.Cells.Item($r, 7).Value = $categoria
    $ws.Cells.Item($r, 8).Value = $variedad
    $ws.Cells.Item($r, 9).Value = $calidad
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $precioMin
    $ws.Cells.Item($r, 12).Value = $precioMax
    $ws.Cells.Item($r, 13).Value = $precioProm
    $ws.Cells.Item($r, 14).Value = $unidadCom
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $kgUnid
    $ws.Cells.Item($r, 18).Value = $tipo
}

Set-Row 367 44706 "Primera" 500  700 800 750
Set-Row 368 44706 "Segunda" 800  600 700 650
Set-Row 369 44706 "Tercera" 1200 400 500 450

$ws.Range("D367:D369").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "done"
